$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update "days" column (D) for GERPHIS (row 7) and STSWENG (row 8)
$ws.Range("D7").Value = "T Th"
$ws.Range("D8").Value = "T F"

# Update the active selection to match the saved state
$ws.Range("F7").Select()
